$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.543.75'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.849.12'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.50'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6296'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07456'
$ws.Range('E8').Value = '  -1.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2909'
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.00'
$ws.Range('E10').Value = '  +1.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07746'
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').Value = '1.853.19'
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.013'
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6828'
$ws.Range('E14').Value = '  +0.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001021'
$ws.Range('E15').Value = '  -1.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.66'
$ws.Range('E16').Value = '  -0.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.315'
$ws.Range('E17').Value = '  +3.33%  '
$ws.Range('D18').Value = '29.554.74'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '229.90'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.37'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9998'
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.515'
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '159.39'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.508'
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1364'
$ws.Range('E26').Value = '  -2.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.54'
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06604'
$ws.Range('E28').Value = '  +16.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.471'
$ws.Range('E29').Value = '  +2.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.489'
$ws.Range('E30').Value = '  +0.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.102'
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.092'
$ws.Range('E32').Value = '  +1.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.848'
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('E34').Value = '  -1.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6970'
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.564'
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('E37').Value = '  +2.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.838'
$ws.Range('E38').Value = '  +4.42%  '
$ws.Range('D39').Value = '1.254.15'
$ws.Range('E39').Value = '  +1.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.789'
$ws.Range('E40').Value = '  +5.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9346'
$ws.Range('E41').Value = '  +3.67%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').Value = '2.011.26'
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.26'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.098'
$ws.Range('E46').Value = '  -0.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.726'
$ws.Range('E47').Value = '  +2.64%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1156'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.988'
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3937'
$ws.Range('E50').Value = '  -1.48%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.00000000114'
$ws.Range('E51').Value = '  -1.58%  '
